# Add a new data row (row 4) to the "Tab_6a_Zeitreihen" sheet, matching the
# layout/style of the existing data rows, and widen columns C and D so the
# new (longer) German/English descriptions are legible.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 4: copy the formatting from row 3 (same style used by every
# data row: A/B/C/D use cell style index 4) then fill in the new values. ---
$ws.Range("A3:D3").Copy()
$ws.Range("A4:D4").PasteSpecial(-4122)

$ws.Range("A4").Value = "Z07_B02_P01_Ib01_I01_Z01"
$ws.Range("B4").Value = "Z07_B02_P01_Ib01_I01"
$ws.Range("C4").Value = "Erzeugung erneuerbarer Energien in Relation zum Brutto-Endenergieverbrauch"
$ws.Range("D4").Value = "Generation of renewable energies as a share of gross final energy consumption"

# --- Widen columns C (~27.39 chars) and D (~32.96 chars) to fit the new,
# longer text that was just added. ---
$ws.Columns.Item(3).ColumnWidth = 26.714285714285715
$ws.Columns.Item(4).ColumnWidth = 32.285714285714285
